# Scheduled market-data refresh: update computed price/profit columns (H:N)
# across the Leve-profit worksheets. Values sourced from the latest Universalis
# pull; some leves currently have no listings, so their derived cells are cleared
# instead of left at a stale 0/blank mix.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 697.44446
$ws.Cells.Item(28, 9).Value = 583.2857
$ws.Cells.Item(28, 10).Value = 1097
$ws.Cells.Item(28, 11).Value = 583.2857
$ws.Cells.Item(28, 12).Value = 1097
$ws.Cells.Item(28, 13).Value = -98.28570000000002
$ws.Cells.Item(28, 14).Value = -2067
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = $null
$ws.Cells.Item(43, 14).Value = $null
$ws.Cells.Item(62, 8).Value = 3436
$ws.Cells.Item(62, 9).Value = 2918.25
$ws.Cells.Item(62, 10).Value = 4126.3335
$ws.Cells.Item(62, 11).Value = 2918.25
$ws.Cells.Item(62, 12).Value = 4126.3335
$ws.Cells.Item(62, 13).Value = -2294.25
$ws.Cells.Item(62, 14).Value = -5374.3335
$ws.Cells.Item(65, 8).Value = 3436
$ws.Cells.Item(65, 9).Value = 2918.25
$ws.Cells.Item(65, 10).Value = 4126.3335
$ws.Cells.Item(65, 11).Value = 14591.25
$ws.Cells.Item(65, 12).Value = 20631.6675
$ws.Cells.Item(65, 13).Value = -11471.25
$ws.Cells.Item(65, 14).Value = -26871.6675
$ws.Cells.Item(80, 8).Value = 1478.1333
$ws.Cells.Item(80, 9).Value = 2093.6667
$ws.Cells.Item(80, 10).Value = 554.8333
$ws.Cells.Item(80, 11).Value = 6281.000100000001
$ws.Cells.Item(80, 12).Value = 1664.4999
$ws.Cells.Item(80, 13).Value = -5283.000100000001
$ws.Cells.Item(80, 14).Value = -3660.4999
$ws.Cells.Item(83, 8).Value = 1478.1333
$ws.Cells.Item(83, 9).Value = 2093.6667
$ws.Cells.Item(83, 10).Value = 554.8333
$ws.Cells.Item(83, 11).Value = 18843.0003
$ws.Cells.Item(83, 12).Value = 4993.4997
$ws.Cells.Item(83, 13).Value = -13851.0003
$ws.Cells.Item(83, 14).Value = -14977.4997
$ws.Cells.Item(135, 8).Value = 29412272
$ws.Cells.Item(135, 9).Value = 31250516
$ws.Cells.Item(135, 10).Value = 376
$ws.Cells.Item(135, 11).Value = 281254644
$ws.Cells.Item(135, 12).Value = 3384
$ws.Cells.Item(135, 13).Value = -281252109
$ws.Cells.Item(135, 14).Value = -8454
$ws.Cells.Item(138, 8).Value = 2020.711
$ws.Cells.Item(138, 9).Value = 1690.7826
$ws.Cells.Item(138, 10).Value = 2365.6365
$ws.Cells.Item(138, 11).Value = 5072.3478
$ws.Cells.Item(138, 12).Value = 7096.9095
$ws.Cells.Item(138, 13).Value = 67.65220000000045
$ws.Cells.Item(138, 14).Value = -17376.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2413
$ws.Cells.Item(63, 9).Value = 2411.4285
$ws.Cells.Item(63, 10).Value = 2424
$ws.Cells.Item(63, 11).Value = 2411.4285
$ws.Cells.Item(63, 12).Value = 2424
$ws.Cells.Item(63, 13).Value = -1725.4285
$ws.Cells.Item(63, 14).Value = -3796
$ws.Cells.Item(66, 8).Value = 2413
$ws.Cells.Item(66, 9).Value = 2411.4285
$ws.Cells.Item(66, 10).Value = 2424
$ws.Cells.Item(66, 11).Value = 12057.1425
$ws.Cells.Item(66, 12).Value = 12120
$ws.Cells.Item(66, 13).Value = -8625.1425
$ws.Cells.Item(66, 14).Value = -18984
$ws.Cells.Item(74, 8).Value = 32260428
$ws.Cells.Item(74, 9).Value = 34485216
$ws.Cells.Item(74, 10).Value = 1011
$ws.Cells.Item(74, 11).Value = 34485216
$ws.Cells.Item(74, 12).Value = 1011
$ws.Cells.Item(74, 13).Value = -34484342
$ws.Cells.Item(74, 14).Value = -2759
$ws.Cells.Item(77, 8).Value = 32260428
$ws.Cells.Item(77, 9).Value = 34485216
$ws.Cells.Item(77, 10).Value = 1011
$ws.Cells.Item(77, 11).Value = 172426080
$ws.Cells.Item(77, 12).Value = 5055
$ws.Cells.Item(77, 13).Value = -172421712
$ws.Cells.Item(77, 14).Value = -13791
$ws.Cells.Item(110, 8).Value = 126227.125
$ws.Cells.Item(110, 9).Value = 143545.28
$ws.Cells.Item(110, 10).Value = 5000
$ws.Cells.Item(110, 11).Value = 143545.28
$ws.Cells.Item(110, 12).Value = 5000
$ws.Cells.Item(110, 13).Value = -141500.28
$ws.Cells.Item(110, 14).Value = -9090
$ws.Cells.Item(122, 8).Value = 4019.6924
$ws.Cells.Item(122, 9).Value = 3180.48
$ws.Cells.Item(122, 10).Value = 25000
$ws.Cells.Item(122, 11).Value = 9541.440000000001
$ws.Cells.Item(122, 12).Value = 75000
$ws.Cells.Item(122, 13).Value = -7091.440000000001
$ws.Cells.Item(122, 14).Value = -79900
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = $null
$ws.Cells.Item(124, 14).Value = $null
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = $null
$ws.Cells.Item(129, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 3849328.8
$ws.Cells.Item(132, 9).Value = 4350850
$ws.Cells.Item(132, 10).Value = 4333.3335
$ws.Cells.Item(132, 11).Value = 13052550
$ws.Cells.Item(132, 12).Value = 13000.0005
$ws.Cells.Item(132, 13).Value = -13050020
$ws.Cells.Item(132, 14).Value = -18060.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4236.7646
$ws.Cells.Item(86, 9).Value = 4194.8
$ws.Cells.Item(86, 10).Value = 4296.7144
$ws.Cells.Item(86, 11).Value = 4194.8
$ws.Cells.Item(86, 12).Value = 4296.7144
$ws.Cells.Item(86, 13).Value = -3071.8
$ws.Cells.Item(86, 14).Value = -6542.7144
$ws.Cells.Item(89, 8).Value = 4236.7646
$ws.Cells.Item(89, 9).Value = 4194.8
$ws.Cells.Item(89, 10).Value = 4296.7144
$ws.Cells.Item(89, 11).Value = 20974
$ws.Cells.Item(89, 12).Value = 21483.572
$ws.Cells.Item(89, 13).Value = -15358
$ws.Cells.Item(89, 14).Value = -32715.572
$ws.Cells.Item(134, 8).Value = 20838510
$ws.Cells.Item(134, 9).Value = 22732498
$ws.Cells.Item(134, 10).Value = 4645.5
$ws.Cells.Item(134, 11).Value = 68197494
$ws.Cells.Item(134, 12).Value = 13936.5
$ws.Cells.Item(134, 13).Value = -68194959
$ws.Cells.Item(134, 14).Value = -19006.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 500003740
$ws.Cells.Item(4, 9).Value = 500003740
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 500003740
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -500003628
$ws.Cells.Item(4, 14).Value = $null
$ws.Cells.Item(16, 8).Value = 1858.2142
$ws.Cells.Item(16, 9).Value = 1766.8334
$ws.Cells.Item(16, 10).Value = 2406.5
$ws.Cells.Item(16, 11).Value = 1766.8334
$ws.Cells.Item(16, 12).Value = 2406.5
$ws.Cells.Item(16, 13).Value = -1479.8334
$ws.Cells.Item(16, 14).Value = -2980.5
$ws.Cells.Item(22, 8).Value = 5459.048
$ws.Cells.Item(22, 9).Value = 6906.1875
$ws.Cells.Item(22, 10).Value = 828.2
$ws.Cells.Item(22, 11).Value = 6906.1875
$ws.Cells.Item(22, 12).Value = 828.2
$ws.Cells.Item(22, 13).Value = -6556.1875
$ws.Cells.Item(22, 14).Value = -1528.2
$ws.Cells.Item(31, 8).Value = 6743
$ws.Cells.Item(31, 9).Value = 3991.6
$ws.Cells.Item(31, 10).Value = 20500
$ws.Cells.Item(31, 11).Value = 3991.6
$ws.Cells.Item(31, 12).Value = 20500
$ws.Cells.Item(31, 13).Value = -3696.6
$ws.Cells.Item(31, 14).Value = -21090
$ws.Cells.Item(34, 8).Value = 6743
$ws.Cells.Item(34, 9).Value = 3991.6
$ws.Cells.Item(34, 10).Value = 20500
$ws.Cells.Item(34, 11).Value = 3991.6
$ws.Cells.Item(34, 12).Value = 20500
$ws.Cells.Item(34, 13).Value = -3789.6
$ws.Cells.Item(34, 14).Value = -20904
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = $null
$ws.Cells.Item(48, 14).Value = $null
$ws.Cells.Item(54, 8).Value = 34999.75
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 34999.75
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 34999.75
$ws.Cells.Item(54, 14).Value = -36315.75
$ws.Cells.Item(62, 8).Value = 3229.7778
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 3229.7778
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 3229.7778
$ws.Cells.Item(62, 13).Value = $null
$ws.Cells.Item(62, 14).Value = -4477.7778
$ws.Cells.Item(65, 8).Value = 3229.7778
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 3229.7778
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 16148.889
$ws.Cells.Item(65, 13).Value = $null
$ws.Cells.Item(65, 14).Value = -22388.889
$ws.Cells.Item(99, 8).Value = 2840
$ws.Cells.Item(99, 9).Value = 2840
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 2840
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = $null
$ws.Cells.Item(99, 14).Value = $null
$ws.Cells.Item(113, 8).Value = 1858.2142
$ws.Cells.Item(113, 9).Value = 1766.8334
$ws.Cells.Item(113, 10).Value = 2406.5
$ws.Cells.Item(113, 11).Value = 1766.8334
$ws.Cells.Item(113, 12).Value = 2406.5
$ws.Cells.Item(113, 13).Value = 403.1666
$ws.Cells.Item(113, 14).Value = -6746.5
$ws.Cells.Item(126, 8).Value = 2840
$ws.Cells.Item(126, 9).Value = 2840
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8520
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = $null
$ws.Cells.Item(126, 14).Value = $null
$ws.Cells.Item(134, 8).Value = 50002800
$ws.Cells.Item(134, 9).Value = 50002800
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 150008400
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -150005865

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 6474.231
$ws.Cells.Item(3, 9).Value = 5833.1816
$ws.Cells.Item(3, 10).Value = 10000
$ws.Cells.Item(3, 11).Value = 17499.5448
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 13).Value = -17387.5448
$ws.Cells.Item(3, 14).Value = -30224
$ws.Cells.Item(11, 8).Value = 148809.53
$ws.Cells.Item(11, 9).Value = 151250
$ws.Cells.Item(11, 10).Value = 100000
$ws.Cells.Item(11, 11).Value = 453750
$ws.Cells.Item(11, 12).Value = 300000
$ws.Cells.Item(11, 13).Value = -453610
$ws.Cells.Item(11, 14).Value = -300280
$ws.Cells.Item(56, 8).Value = 12845.109
$ws.Cells.Item(56, 9).Value = 12845.109
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 12845.109
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = -12315.109

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 25000
$ws.Cells.Item(5, 9).Value = 25000
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = $null
$ws.Cells.Item(5, 14).Value = $null
$ws.Cells.Item(70, 8).Value = 6956.5386
$ws.Cells.Item(70, 9).Value = 6873.75
$ws.Cells.Item(70, 10).Value = 7950
$ws.Cells.Item(70, 11).Value = 6873.75
$ws.Cells.Item(70, 12).Value = 7950
$ws.Cells.Item(70, 13).Value = -6603.75
$ws.Cells.Item(70, 14).Value = -8490
$ws.Cells.Item(73, 8).Value = 6956.5386
$ws.Cells.Item(73, 9).Value = 6873.75
$ws.Cells.Item(73, 10).Value = 7950
$ws.Cells.Item(73, 11).Value = 6873.75
$ws.Cells.Item(73, 12).Value = 7950
$ws.Cells.Item(73, 13).Value = -5937.75
$ws.Cells.Item(73, 14).Value = -9822
$ws.Cells.Item(126, 8).Value = 2922.25
$ws.Cells.Item(126, 9).Value = 2922.25
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8766.75
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -6296.75
$ws.Cells.Item(139, 8).Value = 119241.5
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 119241.5
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 119241.5
$ws.Cells.Item(139, 14).Value = -129521.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 866.5
$ws.Cells.Item(2, 9).Value = 499
$ws.Cells.Item(2, 10).Value = 1234
$ws.Cells.Item(2, 11).Value = 499
$ws.Cells.Item(2, 12).Value = 1234
$ws.Cells.Item(2, 13).Value = -387
$ws.Cells.Item(2, 14).Value = -1458

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 10002780
$ws.Cells.Item(132, 9).Value = 11906460
$ws.Cells.Item(132, 10).Value = 8459.125
$ws.Cells.Item(132, 11).Value = 35719380
$ws.Cells.Item(132, 12).Value = 25377.375
$ws.Cells.Item(132, 13).Value = -35716850
$ws.Cells.Item(132, 14).Value = -30437.375
